$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.608.63'
$ws.Range('E2').Value = '  -6.61%  '
$ws.Range('D3').Value = '2.891.63'
$ws.Range('E3').Value = '  -5.01%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '550.40'
$ws.Range('E5').Value = '  -5.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '121.18'
$ws.Range('E6').Value = '  -6.72%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '2.882.62'
$ws.Range('E8').Value = '  -5.22%  '
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.124'
$ws.Range('E10').Value = '  -9.72%  '
$ws.Range('E11').Value = '  -9.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.432'
$ws.Range('E12').Value = '  -1.37%  '
$ws.Range('E13').Value = '  -9.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '31.46'
$ws.Range('E14').Value = '  -5.85%  '
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('D16').Value = '3.366.31'
$ws.Range('E16').Value = '  -5.07%  '
$ws.Range('D17').Value = '2.889.80'
$ws.Range('E17').Value = '  -5.14%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '57.567.61'
$ws.Range('E18').Value = '  -6.72%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.47'
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '408.38'
$ws.Range('E20').Value = '  -8.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.78'
$ws.Range('E21').Value = '  -4.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.650'
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('E23').Value = '  -8.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.55'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '76.63'
$ws.Range('E25').Value = '  -5.02%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('E28').Value = '  -4.14%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.12'
$ws.Range('E29').Value = '  -3.99%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.90'
$ws.Range('E30').Value = '  -4.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.01'
$ws.Range('E31').Value = '  -6.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '24.53'
$ws.Range('E32').Value = '  -5.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0949'
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.897'
$ws.Range('E34').Value = '  -7.51%  '
$ws.Range('B35').Value = 'Stacks'
$ws.Range('C35').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.01'
$ws.Range('E35').Value = '  -13.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.33'
$ws.Range('E36').Value = '  -6.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '48.29'
$ws.Range('E37').Value = '  -3.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.33'
$ws.Range('E38').Value = '  +5.57%  '
$ws.Range('D39').Value = '0.0₃0616'
$ws.Range('E39').Value = '  -11.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0342'
$ws.Range('E40').Value = '  -8.00%  '
$ws.Range('E41').Value = '  -4.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '362.45'
$ws.Range('E42').Value = '  -3.92%  '
$ws.Range('D43').Value = '2.601.66'
$ws.Range('E43').Value = '  -3.07%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.33'
$ws.Range('E45').Value = '  -7.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '117.58'
$ws.Range('E46').Value = '  -4.30%  '
$ws.Range('E47').Value = '  -4.63%  '
$ws.Range('E48').Value = '  -1.90%  '
$ws.Range('E49').Value = '  -3.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.10'
$ws.Range('E50').Value = '  -7.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.94'
$ws.Range('E51').Value = '  -5.31%  '
